# PHOENIX-5876: refactoring the water charge management module
#
# The "budgetCheck" fixture row (row 5) on the journalVoucherDetails sheet
# carries a voucherDate that needs to be bumped from 2016 to 2017 so the
# water-charge functional test data stays valid. Updating that cell and
# re-selecting/activating the journalVoucherDetails sheet (instead of the
# financialBankDetails sheet that used to be active) is the whole of the
# functional change.

$wb = $excel.ActiveWorkbook

$journalVoucherDetails = $wb.Worksheets.Item("journalVoucherDetails")

# Bump the budgetCheck row's voucherDate from 03/01/2016 to 03/01/2017.
$journalVoucherDetails.Range("B5").Value = "03/01/2017"

# journalVoucherDetails becomes the active sheet/tab with B5 selected and
# scrolled back to the top-left (it previously had F1 scrolled into view
# with H6 selected, while financialBankDetails used to be the active tab).
$journalVoucherDetails.Activate()
$journalVoucherDetails.Range("B5").Select()
